# Add uploaded rows of data (Group Number / Group Name pairs) to Sheet1.
# Rows 3-12 are appended below the existing header/data rows (1-2), growing
# the used range from A1:B2 to A1:B12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append, starting at row 3. Column A values are numeric-looking
# identifiers that must be kept as text (same as existing rows 1-2 which
# are flagged via ignoredErrors numberStoredAsText).
$data = @(
    @("6314",   "Bhaiya"),
    @("145",    "SIS"),
    @("3456",   "vinag"),
    @("450",    "vikram"),
    @("8614",   "jaiin"),
    @("11858",  "Dipanshu"),
    @("945087", "Bai"),
    @("731713", "Vi"),
    @("0000",   "vuh"),
    @("888",    "Dipanshu madd")
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $numberCell = $ws.Cells.Item($row, 1)
    $numberCell.NumberFormat = "@"
    $numberCell.Value = $data[$i][0]

    $nameCell = $ws.Cells.Item($row, 2)
    $nameCell.Value = $data[$i][1]
}

Write-Host "Added $($data.Length) rows of data"
